# Edit script for Brentford_stats.xlsx
#
# 1. Renames the per-statistic worksheet tabs to their more readable,
#    spaced-out forms (the "Matches" and "Possession" sheets are unchanged).
# 2. Bumps the day component of every player's "Age" column (format
#    "YY-DDD", e.g. "23-324") forward by one day on every statistics
#    sheet, leaving summary rows (e.g. "Squad Total" / "Opponent Total",
#    whose Age is a plain decimal average) untouched.

$wb = $excel.ActiveWorkbook

# --- 1. Rename worksheet tabs -------------------------------------------------

$renames = @{
    "StandardStats"    = "Standard Stats";
    "ShootingStats"    = "Shooting Stats";
    "PassingStats"     = "Passing Stats";
    "PassTypes"        = "Pass Types";
    "GoalShotCreation" = "Goal & Shot Creation";
    "DefensiveActions" = "Defensive Actions";
    "PlayingTime"      = "Playing Time";
    "MiscStats"        = "Miscellaneous Stats";
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}

# --- 2. Increment the "day" portion of every Age value ("YY-DDD") ------------

$ageColumn = 5   # column E
$firstDataRow = 4

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Matches") {
        continue
    }

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

    for ($row = $firstDataRow; $row -le $lastRow; $row++) {
        $cell = $ws.Cells.Item($row, $ageColumn)
        $text = [string]$cell.Value2

        if ($text -match '^(\d+)-(\d+)$') {
            $years = $matches[1]
            $newDay = [int]$matches[2] + 1
            $newDayText = $newDay.ToString().PadLeft(3, '0')
            $cell.Value2 = $years + "-" + $newDayText
        }
    }
}
